# Apply the cryptocurrency price / 1h-volume updates described in the commit diff.
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.448.75"
$ws.Cells.Item(2, 5).Value = "  +1.08%  "

$ws.Cells.Item(3, 4).Value = "'3.094.98"
$ws.Cells.Item(3, 5).Value = "  -0.22%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).Value = "'582.88"
$ws.Cells.Item(5, 5).Value = "  -0.23%  "

$ws.Cells.Item(6, 4).Value = "'145.02"
$ws.Cells.Item(6, 5).Value = "  +1.61%  "

$ws.Cells.Item(7, 5).Value = "  +0.09%  "

$ws.Cells.Item(8, 4).Value = "'3.086.89"
$ws.Cells.Item(8, 5).Value = "  -0.34%  "

$ws.Cells.Item(9, 4).Value = "'0.527"
$ws.Cells.Item(9, 5).Value = "  -0.22%  "

$ws.Cells.Item(10, 4).Value = "'0.158"
$ws.Cells.Item(10, 5).Value = "  +6.24%  "

$ws.Cells.Item(11, 4).Value = "'5.65"
$ws.Cells.Item(11, 5).Value = "  -0.68%  "

$ws.Cells.Item(12, 4).Value = "'0.456"
$ws.Cells.Item(12, 5).Value = "  -2.32%  "

$ws.Cells.Item(13, 4).Value = "'0.0000246"
$ws.Cells.Item(13, 5).Value = "  +0.80%  "

$ws.Cells.Item(14, 4).Value = "'37.42"
$ws.Cells.Item(14, 5).Value = "  +6.00%  "

$ws.Cells.Item(15, 5).Value = "  -1.03%  "

$ws.Cells.Item(16, 4).Value = "'3.607.69"
$ws.Cells.Item(16, 5).Value = "  -0.20%  "

$ws.Cells.Item(17, 4).Value = "'63.303.40"
$ws.Cells.Item(17, 5).Value = "  +0.97%  "

$ws.Cells.Item(18, 4).Value = "'7.09"
$ws.Cells.Item(18, 5).Value = "  -1.08%  "

$ws.Cells.Item(19, 4).Value = "'3.092.96"
$ws.Cells.Item(19, 5).Value = "  -0.29%  "

$ws.Cells.Item(20, 4).Value = "'459.67"
$ws.Cells.Item(20, 5).Value = "  -0.40%  "

$ws.Cells.Item(21, 4).Value = "'14.24"
$ws.Cells.Item(21, 5).Value = "  +1.56%  "

$ws.Cells.Item(22, 4).Value = "'0.725"
$ws.Cells.Item(22, 5).Value = "  -0.32%  "

$ws.Cells.Item(23, 4).Value = "'7.44"
$ws.Cells.Item(23, 5).Value = "  -1.12%  "

$ws.Cells.Item(24, 4).Value = "'12.97"
$ws.Cells.Item(24, 5).Value = "  -3.25%  "

$ws.Cells.Item(25, 4).Value = "'81.19"
$ws.Cells.Item(25, 5).Value = "  -0.93%  "

$ws.Cells.Item(26, 4).Value = "'2.12"
$ws.Cells.Item(26, 5).Value = "  -2.18%  "

$ws.Cells.Item(27, 5).Value = "  +0.06%  "

$ws.Cells.Item(28, 4).Value = "'8.86"
$ws.Cells.Item(28, 5).Value = "  +7.25%  "

$ws.Cells.Item(29, 5).Value = "  +0.07%  "

$ws.Cells.Item(30, 4).Value = "'2.67"
$ws.Cells.Item(30, 5).Value = "  -0.80%  "

$ws.Cells.Item(31, 4).Value = "'2.20"
$ws.Cells.Item(31, 5).Value = "  -1.74%  "

$ws.Cells.Item(32, 4).Value = "'6.78"
$ws.Cells.Item(32, 5).Value = "  -0.44%  "

$ws.Cells.Item(33, 4).Value = "'26.71"
$ws.Cells.Item(33, 5).Value = "  -0.67%  "

$ws.Cells.Item(34, 4).Value = "'0.107"
$ws.Cells.Item(34, 5).Value = "  -2.61%  "

$ws.Cells.Item(35, 4).Value = "'0.0₃0847"
$ws.Cells.Item(35, 5).Value = "  +3.67%  "

$ws.Cells.Item(36, 2).Value = "Stacks"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(36, 4).Value = "'2.31"
$ws.Cells.Item(36, 5).Value = "  -2.33%  "

$ws.Cells.Item(37, 2).Value = "Mantle"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(37, 4).Value = "'1.03"
$ws.Cells.Item(37, 5).Value = "  -1.00%  "

$ws.Cells.Item(38, 4).Value = "'3.36"
$ws.Cells.Item(38, 5).Value = "  +6.94%  "

$ws.Cells.Item(39, 4).Value = "'6.02"
$ws.Cells.Item(39, 5).Value = "  -0.23%  "

$ws.Cells.Item(40, 4).Value = "'50.19"
$ws.Cells.Item(40, 5).Value = "  -1.55%  "

$ws.Cells.Item(41, 4).Value = "'438.23"
$ws.Cells.Item(41, 5).Value = "  +2.96%  "

$ws.Cells.Item(42, 4).Value = "'8.72"
$ws.Cells.Item(42, 5).Value = "  -0.62%  "

$ws.Cells.Item(43, 4).Value = "'0.0368"
$ws.Cells.Item(43, 5).Value = "  +0.09%  "

$ws.Cells.Item(44, 4).Value = "'2.859.82"
$ws.Cells.Item(44, 5).Value = "  -1.76%  "

$ws.Cells.Item(45, 2).Value = "Kaspa"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(45, 4).Value = "'0.108"
$ws.Cells.Item(45, 5).Value = "  -1.54%  "

$ws.Cells.Item(46, 2).Value = "TheGraph"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(46, 4).Value = "'0.269"
$ws.Cells.Item(46, 5).Value = "  -2.83%  "

$ws.Cells.Item(47, 4).Value = "'35.84"
$ws.Cells.Item(47, 5).Value = "  +2.77%  "

$ws.Cells.Item(48, 5).Value = "  +0.05%  "

$ws.Cells.Item(49, 4).Value = "'123.55"
$ws.Cells.Item(49, 5).Value = "  -0.26%  "

$ws.Cells.Item(50, 5).Value = "  -1.12%  "

$ws.Cells.Item(51, 4).Value = "'24.12"
$ws.Cells.Item(51, 5).Value = "  -2.47%  "
